# Rename the embedded logo picture objects:
#   image1.png -> image2.png  (Pearson Edexcel logo, appears in both footers)
#   image2.jpg -> image1.jpg  (BTEC logo, appears in the first-page header)
#
# InlineShape has no writable .Name property in the Word object model, so we
# round-trip each picture through Shape (InlineShape.ConvertToShape /
# Shape.ConvertToInlineShape), which does expose a settable .Name.

$d = $word.ActiveDocument

$renameMap = @{
    "image1.png" = "image2.png"
    "image2.jpg" = "image1.jpg"
}

function Rename-ShapesInRange($range, [string]$label) {
    $count = $range.InlineShapes.Count
    for ($k = 1; $k -le $count; $k++) {
        $inlineShape = $range.InlineShapes.Item($k)
        $shape = $inlineShape.ConvertToShape()
        $currentName = $shape.Name
        if ($renameMap.ContainsKey($currentName)) {
            $shape.Name = $renameMap[$currentName]
        }
        [void]$shape.ConvertToInlineShape()
    }
}

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $section = $d.Sections.Item($si)

    for ($hi = 1; $hi -le 3; $hi++) {
        $header = $section.Headers.Item($hi)
        if ($header.Exists) {
            Rename-ShapesInRange $header.Range "header $hi"
        }
    }

    for ($fi = 1; $fi -le 3; $fi++) {
        $footer = $section.Footers.Item($fi)
        if ($footer.Exists) {
            Rename-ShapesInRange $footer.Range "footer $fi"
        }
    }
}

Write-Host "Done renaming logo picture objects."
